# chore: description of level 2
#
# Widens a couple of boxes/connectors around the "Woopec.Core" box on
# slide 2 to make room for its renamed label "Woopec.Graphics", and
# nudges a couple of connectors/labels to match.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$emuPerPt = 12700.0

# Shape 5 ("Rechteck 1") - dashed outer rectangle: widen it.
$shp = $s.Shapes.Item(5)
$shp.Width = 8087646 / $emuPerPt

# Shape 6 ("Rechteck 4") - the "Woopec.Core" -> "Woopec.Graphics" box:
# move it left a bit and widen it, then rename the text.
$shp = $s.Shapes.Item(6)
$shp.Left = 6667498 / $emuPerPt
$shp.Width = 2012816 / $emuPerPt
$shp.TextFrame.TextRange.Text = "Woopec.Graphics"

# Shape 10 ("Gerade Verbindung mit Pfeil 19") - connector into the box above.
$shp = $s.Shapes.Item(10)
$shp.Width = 3006658 / $emuPerPt

# Shape 13 ("Gerade Verbindung mit Pfeil 31") - connector under the box.
$shp = $s.Shapes.Item(13)
$shp.Width = 1317489 / $emuPerPt

# Shape 14 ("Textfeld 38") - label next to that connector: move up slightly.
$shp = $s.Shapes.Item(14)
$shp.Top = 3965788 / $emuPerPt

# Shape 21 ("Gerade Verbindung mit Pfeil 3") - another connector further down.
$shp = $s.Shapes.Item(21)
$shp.Width = 1283073 / $emuPerPt
